$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.66"
$ws.Range("E2").Value = "'0.74%"
$ws.Range("D3").Value = "'31.58"
$ws.Range("E3").Value = "'0.83%"
$ws.Range("D4").Value = "'5.095"
$ws.Range("E4").Value = "'-0.88%"
$ws.Range("D5").Value = "'0.07808"
$ws.Range("E5").Value = "'-1.78%"
$ws.Range("D6").Value = "'2.260"
$ws.Range("E6").Value = "'-13.96%"
$ws.Range("D7").Value = "'7.801"
$ws.Range("E7").Value = "'-0.37%"
$ws.Range("D8").Value = "'3.832"
$ws.Range("E8").Value = "'0.05%"
$ws.Range("D9").Value = "'0.9143"
$ws.Range("E9").Value = "'0.55%"
$ws.Range("D10").Value = "'0.1748"
$ws.Range("E10").Value = "'0.91%"
$ws.Range("D11").Value = "'0.07548"
$ws.Range("E11").Value = "'5.77%"
$ws.Range("D12").Value = "'0.09120"
$ws.Range("E12").Value = "'13.57%"
$ws.Range("D13").Value = "'0.03097"
$ws.Range("E13").Value = "'2.64%"
$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'0.53%"
$ws.Range("D15").Value = "'0.001517"
$ws.Range("E15").Value = "'1.66%"
$ws.Range("D16").Value = "'0.005892"
$ws.Range("E16").Value = "'-1.84%"
$ws.Range("D17").Value = "'3.477"
$ws.Range("E17").Value = "'-0.70%"
$ws.Range("D18").Value = "'2.251"
$ws.Range("E18").Value = "'-0.13%"
$ws.Range("E19").Value = "'0.23%"
$ws.Range("E20").Value = "'0.80%"
$ws.Range("D21").Value = "'4.035"
$ws.Range("E21").Value = "'-12.93%"
$ws.Range("D22").Value = "'0.1819"
$ws.Range("E22").Value = "'13.73%"
$ws.Range("D23").Value = "'0.04599"
$ws.Range("E23").Value = "'0.14%"
$ws.Range("D24").Value = "'0.001252"
$ws.Range("E24").Value = "'-0.50%"
$ws.Range("D25").Value = "'0.004459"
$ws.Range("E25").Value = "'0.10%"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("E26").Value = "'5.85%"
$ws.Range("E27").Value = "'-1.34%"
$ws.Range("D39").Value = "'0.01780"
$ws.Range("E39").Value = "'-4.15%"
$ws.Range("D40").Value = "'0.04806"
$ws.Range("E40").Value = "'6.40%"
$ws.Range("D41").Value = "'0.007402"
$ws.Range("E41").Value = "'4.79%"
$ws.Range("E42").Value = "'0.94%"
$ws.Range("D43").Value = "'0.002190"
$ws.Range("E43").Value = "'-2.31%"
$ws.Range("D44").Value = "'0.01021"
$ws.Range("E44").Value = "'-2.14%"
$ws.Range("D45").Value = "'0.00006210"
$ws.Range("E45").Value = "'-3.59%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("E47").Value = "'28.90%"
$ws.Range("D48").Value = "'0.7426"
$ws.Range("E48").Value = "'-9.50%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'-0.04%"
